$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57-69 down to 58-70.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new record's data.
$ws.Range("A57").Value2 = 6
$ws.Range("B57").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C57").Value2 = "Metropolitana"
$ws.Range("D57").Value2 = 44641
$ws.Range("E57").Value2 = 13
$ws.Range("F57").Value2 = 100114007
$ws.Range("G57").Value2 = "Jengibre"
$ws.Range("H57").Value2 = "Sin especificar"
$ws.Range("I57").Value2 = "Primera"
$ws.Range("J57").Value2 = 290
$ws.Range("K57").Value2 = 14000
$ws.Range("L57").Value2 = 14000
$ws.Range("M57").Value2 = 14000
$ws.Range("N57").Value2 = "$/caja 13 kilos"
$ws.Range("O57").Value2 = "Perú"
$ws.Range("P57").Value2 = 1077
$ws.Range("Q57").Value2 = 13
$ws.Range("R57").Value2 = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D.
$ws.Range("D57").NumberFormat = $ws.Range("D58").NumberFormat
